$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.555373
$ws.Range("N2").Value = 3.110746
$ws.Range("O2").Value = 0.3885765569564089
$ws.Range("P2").Value = 0.3526211663203079
$ws.Range("Q2").Value = 0.036952551734
$ws.Range("R2").Value = 0.221715310404
$ws.Range("S2").Value = 0.08499212241052816
$ws.Range("T2").Value = 0.1042864475836364
$ws.Range("O3").Value = 0.0380537990759009
$ws.Range("P3").Value = 0.05179896254485632
$ws.Range("S3").Value = 0.008323387222784617
$ws.Range("T3").Value = 0.01531935773649502
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.289563
$ws.Range("N4").Value = 0.868689
$ws.Range("O4").Value = 0.07234109989177429
$ws.Range("P4").Value = 0.09847095466798705
$ws.Range("Q4").Value = 0.006879437754000001
$ws.Range("R4").Value = 0.061914939786
$ws.Range("S4").Value = 0.01582294018319707
$ws.Range("T4").Value = 0.02912243232490907
$ws.Range("M5").Value = 1.631084
$ws.Range("N5").Value = 3.262168
$ws.Range("O5").Value = 0.4074913251205256
$ws.Range("P5").Value = 0.3697857314267338
$ws.Range("Q5").Value = 0.038751293672
$ws.Range("R5").Value = 0.232507762032
$ws.Range("S5").Value = 0.08912928988085425
$ws.Range("T5").Value = 0.1093628062660905
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04059466666666667
$ws.Range("N6").Value = 0.121784
$ws.Range("O6").Value = 0.01014170607573002
$ws.Range("P6").Value = 0.01380492528774525
$ws.Range("Q6").Value = 0.0009644480906666668
$ws.Range("R6").Value = 0.008680032816000001
$ws.Range("S6").Value = 0.002218263322397856
$ws.Range("T6").Value = 0.00408275723332139
$ws.Range("M7").Value = 0.333811
$ws.Range("N7").Value = 1.001433
$ws.Range("O7").Value = 0.08339551287966027
$ws.Range("P7").Value = 0.1135182597523697
$ws.Range("Q7").Value = 0.007930681738000002
$ws.Range("R7").Value = 0.071376135642
$ws.Range("S7").Value = 0.01824083700435896
$ws.Range("T7").Value = 0.03357261893546558
$ws.Range("M8").Value = 1.555373
$ws.Range("N8").Value = 3.110746
$ws.Range("O8").Value = 0.3885765569564089
$ws.Range("P8").Value = 0.3526211663203079
$ws.Range("Q8").Value = 0.1319912858395
$ws.Range("R8").Value = 0.527965143358
$ws.Range("S8").Value = 0.3035844345458807
$ws.Range("T8").Value = 0.2483347187366715
$ws.Range("O9").Value = 0.0380537990759009
$ws.Range("P9").Value = 0.05179896254485632
$ws.Range("S9").Value = 0.02973041185311629
$ws.Range("T9").Value = 0.0364796048083613
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.289563
$ws.Range("N10").Value = 0.868689
$ws.Range("O10").Value = 0.07234109989177429
$ws.Range("P10").Value = 0.09847095466798705
$ws.Range("Q10").Value = 0.0245727505245
$ws.Range("R10").Value = 0.147436503147
$ws.Range("S10").Value = 0.05651815970857722
$ws.Range("T10").Value = 0.06934852234307799
$ws.Range("M11").Value = 1.631084
$ws.Range("N11").Value = 3.262168
$ws.Range("O11").Value = 0.4074913251205256
$ws.Range("P11").Value = 0.3697857314267338
$ws.Range("Q11").Value = 0.138416234866
$ws.Range("R11").Value = 0.5536649394640001
$ws.Range("S11").Value = 0.3183620352396714
$ws.Range("T11").Value = 0.2604229251606432
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.04059466666666667
$ws.Range("N12").Value = 0.121784
$ws.Range("O12").Value = 0.01014170607573002
$ws.Range("P12").Value = 0.01380492528774525
$ws.Range("Q12").Value = 0.003444924305333334
$ws.Range("R12").Value = 0.020669545832
$ws.Range("S12").Value = 0.007923442753332167
$ws.Range("T12").Value = 0.009722168054423861
$ws.Range("M13").Value = 0.333811
$ws.Range("N13").Value = 1.001433
$ws.Range("O13").Value = 0.08339551287966027
$ws.Range("P13").Value = 0.1135182597523697
$ws.Range("Q13").Value = 0.02832770217650001
$ws.Range("R13").Value = 0.169966213059
$ws.Range("S13").Value = 0.0651546758753013
$ws.Range("T13").Value = 0.07994564081690411